$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Submission sheet: update reportStatus values for a few submissions.
#    (Order below matters - it reproduces the original author's shared
#    string allocation order: F5 first, then F11, then F8, then F14.)
# ---------------------------------------------------------------------------
$sub = $wb.Worksheets.Item("Submission")
$sub.Range("F5").Value = "PENDING_MARKING"
$sub.Range("F11").Value = "MARKED_2"
$sub.Range("F8").Value = "MARKED_1"
$sub.Range("F14").Value = "OVERDUE"

# ---------------------------------------------------------------------------
# 2. Add the new "Report" worksheet as the last tab in the workbook.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$rep = $wb.Worksheets.Add($null, $lastSheet)
$rep.Name = "Report"

# Header row
$rep.Range("A1").Value = "id"

# Column B top-to-bottom
$rep.Range("B1").Value = "reportName"
$rep.Range("B2").Value = "Report"
$rep.Range("B3").Value = "InvestigateReport"
$rep.Range("B4").Value = "FinalYearProject"
$rep.Range("B5").Value = "Capstone1"
$rep.Range("B6").Value = "Capstone2"

# Column A top-to-bottom (ids stored as text, like the other sheets)
$rep.Range("A2").NumberFormat = "@"
$rep.Range("A2").Value = "49053257"
$rep.Range("A3").NumberFormat = "@"
$rep.Range("A3").Value = "26662640"
$rep.Range("A4").NumberFormat = "@"
$rep.Range("A4").Value = "34685929"
$rep.Range("A5").NumberFormat = "@"
$rep.Range("A5").Value = "18449474"
$rep.Range("A6").NumberFormat = "@"
$rep.Range("A6").Value = "63860114"

# Column C top-to-bottom
$rep.Range("C1").Value = "reportPath"
$rep.Range("C2").Value = "src/main/resources/Data/SampleDataXlsx/Report.pdf"
$rep.Range("C3").Value = "src/main/resources/Data/SampleDataXlsx/InvestigateReport.pdf"
$rep.Range("C4").Value = "src/main/resources/Data/SampleDataXlsx/FinalYearProject.pdf"
$rep.Range("C5").Value = "src/main/resources/Data/SampleDataXlsx/Capstone1.pdf"
$rep.Range("C6").Value = "src/main/resources/Data/SampleDataXlsx/Capstone2.pdf"

# Column D (reportType) - D1 header, D2 new value, D3-D6 reuse existing values
$rep.Range("D1").Value = "reportType"
$rep.Range("D2").Value = "REPORT"
$rep.Range("D3").Value = "INVESTIGATION"
$rep.Range("D4").Value = "FINAL_YEAR"
$rep.Range("D5").Value = "CAPSTONE_1"
$rep.Range("D6").Value = "CAPSTONE_2"

# Column widths to match the bestFit layout used elsewhere in the workbook
$rep.Columns.Item(1).ColumnWidth = 8.166666666666666
$rep.Columns.Item(2).ColumnWidth = 16.166666666666668
$rep.Columns.Item(3).ColumnWidth = 59.166666666666664
$rep.Columns.Item(4).ColumnWidth = 14.166666666666666

# ---------------------------------------------------------------------------
# 3. Update the saved selections on the Submission and Consultation sheets.
# ---------------------------------------------------------------------------
$sub.Activate()
$sub.Range("I17").Select()

$con = $wb.Worksheets.Item("Consultation")
$con.Activate()
$con.Range("H17").Select()

# ---------------------------------------------------------------------------
# 4. Leave "Report" as the active/selected tab, matching the author's final
#    view state (new sheet selected at C12).
# ---------------------------------------------------------------------------
$rep.Activate()
$rep.Range("C12").Select()
